$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 59335.875
$ws.Range("J3").Value = 59335.875
$ws.Range("L3").Value = 59335.875
$ws.Range("N3").Value = -59563.875
$ws.Range("H33").Value = 1461.9231
$ws.Range("J33").Value = 1861.625
$ws.Range("L33").Value = 1861.625
$ws.Range("N33").Value = -2319.625
$ws.Range("H98").Value = 111113110
$ws.Range("I98").Value = 111113110
$ws.Range("K98").Value = 111113110
$ws.Range("M98").Value = -111111612
$ws.Range("H102").Value = 59335.875
$ws.Range("J102").Value = 59335.875
$ws.Range("L102").Value = 59335.875
$ws.Range("N102").Value = -65825.875
$ws.Range("H122").Value = 111113110
$ws.Range("I122").Value = 111113110
$ws.Range("K122").Value = 333339330
$ws.Range("M122").Value = -333336880
$ws.Range("H135").Value = 7179.28
$ws.Range("I135").Value = 2645.5454
$ws.Range("K135").Value = 23809.9086
$ws.Range("M135").Value = -21274.9086
$ws.Range("H137").Value = 5822.2
$ws.Range("I137").Value = 3410.75
$ws.Range("J137").Value = 8578.143
$ws.Range("K137").Value = 10232.25
$ws.Range("L137").Value = 25734.429
$ws.Range("M137").Value = -7682.25
$ws.Range("N137").Value = -30834.429

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12501033
$ws.Range("I32").Value = 12821341
$ws.Range("J32").Value = 9000
$ws.Range("K32").Value = 12821341
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = -12821054
$ws.Range("N32").Value = -9574
$ws.Range("H61").Value = 45551596
$ws.Range("I61").Value = 125000860
$ws.Range("K61").Value = 125000860
$ws.Range("M61").Value = -125000648
$ws.Range("H103").Value = 49996.668
$ws.Range("J103").Value = 49996.668
$ws.Range("L103").Value = 49996.668
$ws.Range("N103").Value = -52340.668
$ws.Range("H122").Value = 3749.875
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H136").Value = 45551596
$ws.Range("I136").Value = 125000860
$ws.Range("K136").Value = 375002580
$ws.Range("M136").Value = -375000030

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 33291
$ws.Range("J103").Value = 33291
$ws.Range("L103").Value = 33291
$ws.Range("N103").Value = -35635
$ws.Range("H107").Value = 1472.3077
$ws.Range("I107").Value = 1132.7
$ws.Range("K107").Value = 1132.7
$ws.Range("M107").Value = 787.3
$ws.Range("H134").Value = 36279.273
$ws.Range("I134").Value = 4810.8213
$ws.Range("K134").Value = 14432.4639
$ws.Range("M134").Value = -11897.4639

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1822.25
$ws.Range("I58").Value = 1663.6666
$ws.Range("K58").Value = 1663.6666
$ws.Range("M58").Value = -1460.6666
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H132").Value = 1957.6666
$ws.Range("I132").Value = 1778.7059
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5336.1177
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2806.1177
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 1822.25
$ws.Range("I136").Value = 1663.6666
$ws.Range("K136").Value = 4990.9998
$ws.Range("M136").Value = -2440.9998

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 635.3
$ws.Range("J46").Value = 102
$ws.Range("L46").Value = 306
$ws.Range("N46").Value = -488
$ws.Range("H86").Value = 96
$ws.Range("I86").Value = 96
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 288
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 898
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 96
$ws.Range("I89").Value = 96
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 864
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 5064
$ws.Range("N89").ClearContents()
$ws.Range("H102").Value = 5500
$ws.Range("I102").Value = 5500
$ws.Range("K102").Value = 16500
$ws.Range("M102").Value = -14066
$ws.Range("H103").Value = 8065.6
$ws.Range("J103").Value = 13352
$ws.Range("L103").Value = 40056
$ws.Range("N103").Value = -41814
$ws.Range("H107").Value = 1048.6
$ws.Range("I107").Value = 490
$ws.Range("J107").Value = 1188.25
$ws.Range("K107").Value = 1470
$ws.Range("L107").Value = 3564.75
$ws.Range("M107").Value = 450
$ws.Range("N107").Value = -7404.75
$ws.Range("H118").Value = 5443.2
$ws.Range("I118").Value = 3080
$ws.Range("J118").Value = 7806.4
$ws.Range("K118").Value = 9240
$ws.Range("L118").Value = 23419.2
$ws.Range("M118").Value = -7997
$ws.Range("N118").Value = -25905.2
$ws.Range("H131").Value = 4022.5
$ws.Range("I131").Value = 1597.2
$ws.Range("K131").Value = 4791.6
$ws.Range("M131").Value = 248.3999999999996
$ws.Range("H132").Value = 2661.5833
$ws.Range("I132").Value = 2420.5715
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 21785.1435
$ws.Range("L132").Value = 26991
$ws.Range("M132").Value = -19255.1435
$ws.Range("N132").Value = -32051

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744
$ws.Range("H122").Value = 3814.6667
$ws.Range("I122").Value = 3222
$ws.Range("K122").Value = 9666
$ws.Range("M122").Value = -7216

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 67834.2
$ws.Range("J81").Value = 67834.2
$ws.Range("L81").Value = 67834.2
$ws.Range("N81").Value = -69830.2
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 67834.2
$ws.Range("J84").Value = 67834.2
$ws.Range("L84").Value = 203502.6
$ws.Range("N84").Value = -213486.6
$ws.Range("H132").Value = 52524.25
$ws.Range("I132").Value = 52299.7
$ws.Range("J132").Value = 52748.8
$ws.Range("K132").Value = 156899.1
$ws.Range("L132").Value = 158246.4
$ws.Range("M132").Value = -154369.1
$ws.Range("N132").Value = -163306.4
$ws.Range("H136").Value = 131780.94
$ws.Range("I136").Value = 146042.28
$ws.Range("K136").Value = 438126.84
$ws.Range("M136").Value = -435576.84
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H98").Value = 45663.668
$ws.Range("J98").Value = 45663.668
$ws.Range("L98").Value = 45663.668
$ws.Range("N98").Value = -51653.668
$ws.Range("H113").Value = 1135
$ws.Range("I113").Value = 1156.4615
$ws.Range("J113").Value = 1065.25
$ws.Range("K113").Value = 3469.3845
$ws.Range("L113").Value = 3195.75
$ws.Range("M113").Value = -1299.3845
$ws.Range("N113").Value = -7535.75
$ws.Range("H132").Value = 3655.7693
$ws.Range("I132").Value = 3163.8928
$ws.Range("J132").Value = 4907.8184
$ws.Range("K132").Value = 9491.678400000001
$ws.Range("L132").Value = 14723.4552
$ws.Range("M132").Value = -6961.678400000001
$ws.Range("N132").Value = -19783.4552
$ws.Range("H136").Value = 2005.3
$ws.Range("I136").Value = 2006.625
$ws.Range("K136").Value = 6019.875
$ws.Range("M136").Value = -3469.875

Write-Output "Applied all Behemoth_Profits updates."